$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 308-313 with corrected values
$ws.Range("B308").Value = 520
$ws.Range("C308").Value = 4062
$ws.Range("D308").Value = 51
$ws.Range("E308").Value = 1368

$ws.Range("B309").Value = 319
$ws.Range("C309").Value = 2309
$ws.Range("D309").Value = 32
$ws.Range("E309").Value = 538

$ws.Range("B310").Value = 331
$ws.Range("C310").Value = 2416
$ws.Range("D310").Value = 41
$ws.Range("E310").Value = 578

$ws.Range("B311").Value = 1140
$ws.Range("C311").Value = 8764
$ws.Range("D311").Value = 31
$ws.Range("E311").Value = 930

$ws.Range("B312").Value = 290
$ws.Range("C312").Value = 2158
$ws.Range("D312").Value = 12
$ws.Range("E312").Value = 181

$ws.Range("B313").Value = 142
$ws.Range("C313").Value = 1944
$ws.Range("D313").Value = 10
$ws.Range("E313").Value = 44

# Add new row 314 with latest scraped data
# Force the date-like text to stay as text (matches format of column A in all other rows)
$ws.Range("A314").NumberFormat = "@"
$ws.Range("A314").Value = "11.01.2021"
$ws.Range("B314").Value = 78
$ws.Range("C314").Value = 1576
$ws.Range("D314").Value = 12
$ws.Range("E314").Value = 400
